# edit.ps1 - applies the "feat: refactor and change split" edit:
#   1. Removes the last 6 slides (formerly slide6.xml..slide11.xml),
#      leaving only the first 5 slides in the deck.
#   2. Updates the title slide (slide 1): shrinks the first line's font
#      size and swaps both lines of text for the new reference.
#   3. Updates slides 2-5 (verse slides): resizes the title placeholder
#      to fill the whole slide and swaps in the new verse text.

$p = $ppt.ActivePresentation

# --- 1. Drop the trailing six slides (index 11 down to 6) -----------------
for ($i = $p.Slides.Count; $i -ge 6; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- 2. Slide 1: title/reference text + font size --------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

$newLine1 = "Tonon-kiran'i Solomona "
$newLine2 = "3 : 10"

# First paragraph currently reads "Salamo " (7 chars) -> replace it in place
# so the run structure / pPr stay untouched.
$oldLine1Len = 7
$c1 = $tr1.Characters(1, $oldLine1Len)
$c1.Text = $newLine1

# Second paragraph starts right after the (now resized) first paragraph
# plus its paragraph break.
$line1Len = $newLine1.Length
$start2 = $line1Len + 2
$len2 = $tr1.Length - $start2 + 1
$c2 = $tr1.Characters($start2, $len2)
$c2.Text = $newLine2

# Shrink the font used on the first line from 185pt to 100pt.
$cSize = $tr1.Characters(1, $line1Len)
$cSize.Font.Size = 100

# --- 3. Slides 2-5: verse text + placeholder resize -------------------------
$verses = @(
    "10 Volafotsy no nataony andriny,",
    " ary volamena ny fiankinana,",
    " volomparasy ny ondany,",
    " ny ao anatiny voaravaka amin'ny fanomezam-pitiavana avy tamin'ny zanakavavin'i Jerosalema."
)

for ($idx = 0; $idx -lt $verses.Count; $idx++) {
    $slideNum = $idx + 2
    $s = $p.Slides.Item($slideNum)
    $sh = $s.Shapes.Item(1)

    # Resize the title placeholder to span the whole slide height
    # (2084831 EMU -> 6858000 EMU) without disturbing its (absent) offset.
    $targetHeightPt = 6858000 / 12700.0
    $currentHeightPt = $sh.Height
    $scale = $targetHeightPt / $currentHeightPt
    $sh.ScaleHeight($scale, $true)

    # Replace the run text in place, keeping a single clean run.
    $tr = $sh.TextFrame.TextRange
    $c = $tr.Characters(1, $tr.Length)
    $c.Text = $verses[$idx]
}
